$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values per the diff
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 26
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 20
$ws.Range("B4").Value = 11

# Remove rows 5 and 6 (previously held A5/B5 and A6/B6 data)
$ws.Range("A5:B6").Delete()
